$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("B6").Value = 174
$ws.Range("C6").Value = 16
$ws.Range("F6").Value = 17.24

# Row 7
$ws.Range("B7").Value = 110
$ws.Range("F7").Value = 18.18

# Row 9
$ws.Range("B9").Value = 60
$ws.Range("F9").Value = 16.67

# Row 10
$ws.Range("B10").Value = 78

# Row 12
$ws.Range("B12").Value = 87

# Row 13
$ws.Range("B13").Value = 130
$ws.Range("C13").Value = 2

# Row 14
$ws.Range("B14").Value = 114

# Row 15
$ws.Range("B15").Value = 127
$ws.Range("C15").Value = 3

# Row 17
$ws.Range("B17").Value = 47
$ws.Range("F17").Value = 21.28

# Row 19
$ws.Range("B19").Value = 88

# Row 20
$ws.Range("B20").Value = 26

# Column H best-fit width nudges slightly wider after the data edits above
$ws.Columns.Item(8).ColumnWidth = 13.25

# Selection moves to A5
$ws.Range("A5").Select()
